$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feedstock_to_commodity")

# Fix the off-by-one row reference bug in the cumulative formulas:
# G13 incorrectly pulled from Z14 (the row below) instead of Z13 (its own row).
$ws.Range("G13").Formula = "=`$Z13/21+F13"
$ws.Range("G13").Style = "Normal"

# G14 incorrectly pulled from Z15 (the row below) instead of Z14 (its own row).
$ws.Range("G14").Formula = "=`$Z14/21+F14"
$ws.Range("G14").Style = "Normal"

# Select the sheet/cell the author ended up on after making the fix.
$ws.Activate()
$ws.Range("G17").Select()
